$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new columns (Month, Day, Year) between "Transect" (D) and "Date Sampled" (old E).
$ws.Columns("E:G").Insert()

# Make the new columns the same width as the neighbouring "Transect" column (9.5 chars -> 8.67 ColumnWidth).
$ws.Range("E1:G1").ColumnWidth = 8.67

# New header labels for the inserted columns.
$ws.Range("E1").Value = "Month"
$ws.Range("F1").Value = "Day"
$ws.Range("G1").Value = "Year"

# Fill in Month / Day / Year values for each data row, derived from the "Date Sampled" column (now H).
$ws.Range("E2").Value = 8
$ws.Range("F2").Value = 7
$ws.Range("G2").Value = 2014

$ws.Range("E3").Value = 8
$ws.Range("F3").Value = 8
$ws.Range("G3").Value = 2014

$ws.Range("E4").Value = 8
$ws.Range("F4").Value = 8
$ws.Range("G4").Value = 2014

$ws.Range("E5").Value = 8
$ws.Range("F5").Value = 7
$ws.Range("G5").Value = 2014

# Rename the "Start Depth" / "End Depth" headers (now columns M / N) to include units.
$ws.Range("M1").Value = "Start Depth (M)"
$ws.Range("N1").Value = "End Depth (M)"

# Fix data-entry errors: longitudes in this part of the world are negative (west), so correct their sign.
$ws.Range("J2").Value = -70.891210000000001
$ws.Range("L2").Value = -70.8917

$ws.Range("J3").Value = -70.888829999999999
$ws.Range("L3").Value = -70.888440000000003

$ws.Range("J4").Value = -70.889279999999999
$ws.Range("L4").Value = -70.889269999999996

$ws.Range("J5").Value = -70.888170000000002

# Restore the originally selected cell.
$ws.Range("L5").Select() | Out-Null
